# Insert a new data row before the current row 69, shifting the existing
# rows 69-117 down to 70-118 (dimension grows from A1:R117 to A1:R118).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new weekly record.
$ws.Range("A69").Value = 1
$ws.Range("B69").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C69").Value = "Arica y Parinacota"
$ws.Range("D69").Value = 44907
$ws.Range("E69").Value = 15
$ws.Range("F69").Value = 100112021
$ws.Range("G69").Value = "Ají"
$ws.Range("H69").Value = "Inferno"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 160
$ws.Range("K69").Value = 7000
$ws.Range("L69").Value = 8000
$ws.Range("M69").Value = 7500
$ws.Range("N69").Value = "$/caja 15 kilos"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 500
$ws.Range("Q69").Value = 15
$ws.Range("R69").Value = "Hortaliza"

# Keep the date display format consistent with the rest of column D.
$ws.Range("D69").NumberFormat = $ws.Range("D70").NumberFormat
